# Gantt chart update: shift the project start date by one week and
# update the task dates for Phase 2 ("Research and Data Gathering")
# and its sub-tasks, per the refreshed project schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Project Start date (was =DATEVALUE("14/08/2023") -> 8/14/2023); now a
# plain date value one week later (8/21/2023). Entering a literal value
# here replaces the old formula, matching how the date was typed in.
$ws.Range("D3").Value = 45159

# Phase 2 task: "Research and Data Gathering"
$ws.Range("D15").Value = 45194
$ws.Range("E15").Value = 45198

# Sub-tasks under Phase 2
$ws.Range("D16").Value = 45194
$ws.Range("E16").Value = 45198

$ws.Range("D17").Value = 45195
$ws.Range("E17").Value = 45199

$ws.Range("D18").Value = 45197
$ws.Range("E18").Value = 45200

$ws.Range("D19").Value = 45200
$ws.Range("E19").Value = 45206

# Phase 3 sub-tasks
$ws.Range("D21").Value = 45204
$ws.Range("E21").Value = 45207

$ws.Range("D22").Value = 45200
$ws.Range("E22").Value = 45206

$ws.Range("D23").Value = 45200
$ws.Range("E23").Value = 45205

# Restore the user's on-screen selection
$null = $ws.Range("M14").Select()
